$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I9").Value = 99.5
$ws.Range("K9").Value = 99.5
$ws.Range("M9").Value = 69.5
$ws.Range("H62").Value = 6771.778
$ws.Range("I62").Value = 6111.5
$ws.Range("J62").Value = 7300
$ws.Range("K62").Value = 6111.5
$ws.Range("L62").Value = 7300
$ws.Range("M62").Value = -5487.5
$ws.Range("N62").Value = -8548
$ws.Range("H64").Value = 4150
$ws.Range("I64").Value = 3300
$ws.Range("K64").Value = 3300
$ws.Range("M64").Value = -3052
$ws.Range("H65").Value = 6771.778
$ws.Range("I65").Value = 6111.5
$ws.Range("J65").Value = 7300
$ws.Range("K65").Value = 30557.5
$ws.Range("L65").Value = 36500
$ws.Range("M65").Value = -27437.5
$ws.Range("N65").Value = -42740
$ws.Range("H67").Value = 4150
$ws.Range("I67").Value = 3300
$ws.Range("K67").Value = 3300
$ws.Range("M67").Value = -2442
$ws.Range("H116").Value = 3315.4546
$ws.Range("I116").Value = 2421.25
$ws.Range("K116").Value = 2421.25
$ws.Range("M116").Value = 1020.75
$ws.Range("H132").Value = 3966.875
$ws.Range("I132").Value = 3966.875
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11900.625
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9370.625
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 2323.375
$ws.Range("J138").Value = 3000
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 853.75
$ws.Range("H59").Value = 25000
$ws.Range("J59").Value = 25000
$ws.Range("L59").Value = 25000
$ws.Range("N59").Value = -26608
$ws.Range("H97").Value = 385.63635
$ws.Range("I97").Value = 423.1111
$ws.Range("J97").Value = 217
$ws.Range("K97").Value = 423.1111
$ws.Range("L97").Value = 217
$ws.Range("M97").Value = 72.88889999999998
$ws.Range("N97").Value = -1209
$ws.Range("H124").Value = 24500
$ws.Range("J124").Value = 24500
$ws.Range("L124").Value = 24500
$ws.Range("N124").Value = -34320
$ws.Range("H125").Value = 25000
$ws.Range("J125").Value = 25000
$ws.Range("L125").Value = 25000
$ws.Range("N125").Value = -34840
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 250
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H108").Value = 385000
$ws.Range("J108").Value = 385000
$ws.Range("L108").Value = 385000
$ws.Range("N108").Value = -392680
$ws.Range("H134").Value = 4898.7036
$ws.Range("I134").Value = 4939.1665
$ws.Range("J134").Value = 4575
$ws.Range("K134").Value = 14817.4995
$ws.Range("L134").Value = 13725
$ws.Range("M134").Value = -12282.4995
$ws.Range("N134").Value = -18795

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 20129.375
$ws.Range("J60").Value = 19997.5
$ws.Range("L60").Value = 19997.5
$ws.Range("N60").Value = -21019.5
$ws.Range("H92").Value = 62500
$ws.Range("J92").Value = 62500
$ws.Range("L92").Value = 62500
$ws.Range("N92").Value = -67492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 664
$ws.Range("I5").Value = 766.4286
$ws.Range("J5").Value = 425
$ws.Range("K5").Value = 2299.2858
$ws.Range("L5").Value = 1275
$ws.Range("M5").Value = -2187.2858
$ws.Range("N5").Value = -1499
$ws.Range("H23").Value = 179
$ws.Range("I23").Value = 103.666664
$ws.Range("J23").Value = 292
$ws.Range("K23").Value = 310.999992
$ws.Range("L23").Value = 876
$ws.Range("M23").Value = -75.99999200000002
$ws.Range("N23").Value = -1346
$ws.Range("H97").Value = 149.75
$ws.Range("I97").Value = 149.75
$ws.Range("K97").Value = 449.25
$ws.Range("M97").Value = 46.75
$ws.Range("H129").Value = 1020
$ws.Range("J129").Value = 1010
$ws.Range("L129").Value = 3030
$ws.Range("N129").Value = -13030
$ws.Range("H135").Value = 664
$ws.Range("I135").Value = 766.4286
$ws.Range("J135").Value = 425
$ws.Range("K135").Value = 6897.8574
$ws.Range("L135").Value = 3825
$ws.Range("M135").Value = -4362.8574
$ws.Range("N135").Value = -8895

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 5003
$ws.Range("I10").Value = 5003
$ws.Range("K10").Value = 5003
$ws.Range("M10").Value = -4834
$ws.Range("H57").Value = 19997.857
$ws.Range("J57").Value = 19997.5
$ws.Range("L57").Value = 19997.5
$ws.Range("N57").Value = -21637.5
$ws.Range("H80").Value = 2445.625
$ws.Range("I80").Value = 2516.6667
$ws.Range("J80").Value = 2403
$ws.Range("K80").Value = 2516.6667
$ws.Range("L80").Value = 2403
$ws.Range("M80").Value = -1518.6667
$ws.Range("N80").Value = -4399
$ws.Range("H83").Value = 2445.625
$ws.Range("I83").Value = 2516.6667
$ws.Range("J83").Value = 2403
$ws.Range("K83").Value = 12583.3335
$ws.Range("L83").Value = 12015
$ws.Range("M83").Value = -7591.333500000001
$ws.Range("N83").Value = -21999
$ws.Range("H102").Value = 1433
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 8196.6
$ws.Range("J20").Value = 9996
$ws.Range("L20").Value = 9996
$ws.Range("N20").Value = -10448
$ws.Range("H100").Value = 2166.6667
$ws.Range("I100").Value = 2166.6667
$ws.Range("K100").Value = 2166.6667
$ws.Range("M100").Value = -1625.6667
$ws.Range("H132").Value = 3871
$ws.Range("I132").Value = 3032.3
$ws.Range("J132").Value = 6666.6665
$ws.Range("K132").Value = 9096.900000000001
$ws.Range("L132").Value = 19999.9995
$ws.Range("M132").Value = -6566.900000000001
$ws.Range("N132").Value = -25059.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 930
$ws.Range("I96").Value = 930
$ws.Range("K96").Value = 930
$ws.Range("M96").Value = 443
